$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.030.79'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '2.947.46'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '375.92'
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.31'
$ws.Range("E6").Value = '  -2.64%  '
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.591'
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.40'
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0852'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = '3.404.45'
$ws.Range("E13").Value = '  -1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.17'
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.63'
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.35'
$ws.Range("E16").Value = '  +53.09%  '
$ws.Range("D17").Value = '2.944.14'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '50.997.24'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.09'
$ws.Range("E20").Value = '  -5.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.51'
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("D22").Value = '0.0₃0958'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '266.58'
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.90'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.17'
$ws.Range("E25").Value = '  +8.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.20'
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.50'
$ws.Range("E27").Value = '  -3.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.13'
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.72'
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.02'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '50.97'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '33.46'
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0443'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +4.81%  '
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.54'
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("E43").Value = '  -3.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.83'
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.26'
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("E46").Value = '  +3.04%  '
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("E48").Value = '  -1.95%  '
$ws.Range("E49").Value = '  -1.72%  '
$ws.Range("D50").Value = '1.994.07'
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0328'
$ws.Range("E51").Value = '  -1.23%  '
